$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText, $matchCase = $true) {
    $d.Content.Find.Execute($findText, $matchCase, $false, $false, $false, $false,
                             $true, 1, $false, $replaceText, 2)
}

# 1. Update the date/time in the header
Replace-Text "November   8, 2021 (10:39:26 AM)" "November   8, 2021 (07:04:34 PM)"

# 2. "at the screen." -> "on the screen."
Replace-Text "at the screen." "on the screen."

# 3. "Then, remove previous line, and now add this:" -> "Then, remove the previous line, and now add this:"
Replace-Text "Then, remove previous line, and now add this:" "Then, remove the previous line, and now add this:"

# 4. "Remove previous line. Add this line and execute the program:" -> "Remove the previous line. Add this line and execute the program:"
Replace-Text "Remove previous line. Add this line and execute the program:" "Remove the previous line. Add this line and execute the program:"

# 5. "with following values:" -> "with the following values:"
Replace-Text "with following values:" "with the following values:"

# 6. "with value" -> "with the value"
Replace-Text "with value" "with the value"

# 7. "Write a statement to display value stored at index 4. What is that value? Why?"
Replace-Text "Write a statement to display value stored at index 4. What is that value? Why?" "Write a statement to display the value stored at index 4. What is that value? Why?"

# 8. "Write a statement to display characters in the" -> "Write a statement to display the characters in the"
Replace-Text "Write a statement to display characters in the" "Write a statement to display the characters in the"

# 9. "is declared and change" -> "is declared and change the length of"
Replace-Text "is declared and change" "is declared and change the length of"

# 10. "length to" -> "to"
Replace-Text "length to" "to"

# 11. "? (After changing length the first half contains values" -> "? (After changing the length, the first half contains the values"
Replace-Text "? (After changing length the first half contains values" "? (After changing the length, the first half contains the values"

# 12. Final sentence fixes: "you should have" -> "you expected" and "accomodate" -> "accommodate"
Replace-Text "you should have, can you think of a way to perform these array operations in a way that can accomodate arrays of different lengths?" "you expected, can you think of a way to perform these array operations in a way that can accommodate arrays of different lengths?"
